$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = 25

# --- Numeric columns ---
$ws.Cells.Item($r, 1).Value  = 111924482            # A  Id
$ws.Cells.Item($r, 2).Value  = 84741                # B  Taxonsorteringsordning
$ws.Cells.Item($r, 5).Value  = 37                   # E  TaxonId
$ws.Cells.Item($r, 17).Value = 665745.9058803385    # Q  Ost
$ws.Cells.Item($r, 18).Value = 6640602.958293262    # R  Nord
$ws.Cells.Item($r, 19).Value = 10                   # S  Noggrannhet

# --- Plain text columns ---
$ws.Cells.Item($r, 3).Value  = "Ovaliderad"                              # C  Valideringsstatus
$ws.Cells.Item($r, 4).Value  = "NT"                                      # D  Rödlistade
$ws.Cells.Item($r, 6).Value  = "Jättekamskivling"                        # F  Artnamn
$ws.Cells.Item($r, 7).Value  = "Amanita ceciliae"                        # G  Vetenskapligt namn
$ws.Cells.Item($r, 8).Value  = "(Berk. & Broome) Bas"                    # H  Auktor
$ws.Cells.Item($r, 10).Value = "fruktkroppar"                            # J  Enhet
$ws.Cells.Item($r, 16).Value = "Norra Kopphagen (Norra Kopphagen), Upl"  # P  Lokalnamn
$ws.Cells.Item($r, 20).Value = "Uppsala"                                 # T  Län
$ws.Cells.Item($r, 21).Value = "Uppsala"                                 # U  Kommun
$ws.Cells.Item($r, 22).Value = "Uppland"                                 # V  Provins
$ws.Cells.Item($r, 23).Value = "Almunge"                                 # W  Församling
$ws.Cells.Item($r, 49).Value = "Cajsa Björkén"                           # AW Rapportör
$ws.Cells.Item($r, 50).Value = "Cajsa Björkén"                           # AX Artbestämmare

# --- Text columns that look like numbers/dates: force text with a leading
#     apostrophe so Excel doesn't auto-coerce them to a number/date serial ---
$ws.Cells.Item($r, 9).Value  = "'1"            # I  Antal (stored as text "1")
$ws.Cells.Item($r, 25).Value = "'2023-09-06"   # Y  Startdatum
$ws.Cells.Item($r, 26).Value = "'10:18"        # Z  Starttid
$ws.Cells.Item($r, 27).Value = "'2023-09-06"   # AA Slutdatum
$ws.Cells.Item($r, 28).Value = "'10:18"        # AB Sluttid

# --- Boolean columns ---
$ws.Cells.Item($r, 30).Value = $false   # AD Ej återfunnen
$ws.Cells.Item($r, 31).Value = $false   # AE Osäker artbestämning
$ws.Cells.Item($r, 33).Value = $false   # AG Ospontan

# Columns K (11), AT (46) and AY (51) are blank in the source row (empty
# text cells) -- leave them unset so they stay empty.
